$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (zh-cn, de-de) and Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-30-21 06:30:14"

# zh-cn sheet: Status and Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-21 06:30:11"

# de-de sheet: Status and Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-21 06:30:14"
